$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.927.99'
$ws.Range("E2").Value = '  -0.07%  '
$ws.Range("D3").Value = '1.876.46'
$ws.Range("E3").Value = '  -0.86%  '
$ws.Range("D4").Value = '''0.9991'
$ws.Range("D5").Value = '''0.7425'
$ws.Range("E5").Value = '  -4.37%  '
$ws.Range("D6").Value = '''242.77'
$ws.Range("E6").Value = '  -0.48%  '
$ws.Range("D7").Value = '''0.9993'
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("D8").Value = '''0.3162'
$ws.Range("E8").Value = '  +1.04%  '
$ws.Range("D9").Value = '''0.07182'
$ws.Range("E9").Value = '  -0.87%  '
$ws.Range("D10").Value = '''24.83'
$ws.Range("E10").Value = '  -3.75%  '
$ws.Range("D11").Value = '''0.08390'
$ws.Range("E11").Value = '  -3.61%  '
$ws.Range("D12").Value = '''0.7517'
$ws.Range("E12").Value = '  -2.72%  '
$ws.Range("D13").Value = '''5.453'
$ws.Range("E13").Value = '  +1.00%  '
$ws.Range("D14").Value = '1.883.04'
$ws.Range("E14").Value = '  -9.77%  '
$ws.Range("D15").Value = '''92.75'
$ws.Range("E15").Value = '  -1.74%  '
$ws.Range("D16").Value = '29.917.87'
$ws.Range("E16").Value = '  -0.15%  '
$ws.Range("D17").Value = '''6.103'
$ws.Range("E17").Value = '  -1.48%  '
$ws.Range("D18").Value = '''13.59'
$ws.Range("E18").Value = '  -2.28%  '
$ws.Range("D19").Value = '''244.89'
$ws.Range("E19").Value = '  -0.38%  '
$ws.Range("D20").Value = '''0.000007824'
$ws.Range("E20").Value = '  -0.48%  '
$ws.Range("D21").Value = '''0.9986'
$ws.Range("E21").Value = '  -0.23%  '
$ws.Range("D22").Value = '2.122.90'
$ws.Range("E22").Value = '  -7.95%  '
$ws.Range("D23").Value = '''8.017'
$ws.Range("E23").Value = '  -1.60%  '
$ws.Range("D24").Value = '''0.9981'
$ws.Range("E24").Value = '  -0.32%  '
$ws.Range("D25").Value = '''0.1558'
$ws.Range("E25").Value = '  -5.41%  '
$ws.Range("D26").Value = '''9.286'
$ws.Range("E26").Value = '  -2.20%  '
$ws.Range("D27").Value = '''165.24'
$ws.Range("E27").Value = '  +1.42%  '
$ws.Range("D28").Value = '''18.64'
$ws.Range("E28").Value = '  -1.00%  '
$ws.Range("D29").Value = '''2.040'
$ws.Range("E29").Value = '  -0.60%  '
$ws.Range("D30").Value = '''1.518'
$ws.Range("E30").Value = '  +5.81%  '
$ws.Range("D31").Value = '''4.608'
$ws.Range("E31").Value = '  +2.08%  '
$ws.Range("D32").Value = '''1.533'
$ws.Range("E32").Value = '  -0.62%  '
$ws.Range("D33").Value = '''4.277'
$ws.Range("E33").Value = '  +3.71%  '
$ws.Range("D34").Value = '''0.05331'
$ws.Range("E34").Value = '  -2.70%  '
$ws.Range("D35").Value = '''1.239'
$ws.Range("E35").Value = '  -0.47%  '
$ws.Range("D36").Value = '''0.7547'
$ws.Range("E36").Value = '  +0.39%  '
$ws.Range("D37").Value = '''0.9994'
$ws.Range("E37").Value = '  -0.43%  '
$ws.Range("D38").Value = '''2.698'
$ws.Range("E38").Value = '  +0.44%  '
$ws.Range("D39").Value = '''0.01965'
$ws.Range("E39").Value = '  -0.21%  '
$ws.Range("D40").Value = '''2.752'
$ws.Range("E40").Value = '  -1.28%  '
$ws.Range("D41").Value = '''0.4537'
$ws.Range("E41").Value = '  +0.82%  '
$ws.Range("D42").Value = '1.115.28'
$ws.Range("E42").Value = '  +0.80%  '
$ws.Range("D43").Value = '''6.062'
$ws.Range("E43").Value = '  -0.43%  '
$ws.Range("D44").Value = '''72.69'
$ws.Range("E44").Value = '  -1.32%  '
$ws.Range("D45").Value = '''0.8582'
$ws.Range("E45").Value = '  +0.52%  '
$ws.Range("D46").Value = '''1.001'
$ws.Range("E46").Value = '  +0.05%  '
$ws.Range("D47").Value = '''103.57'
$ws.Range("E47").Value = '  +0.35%  '
$ws.Range("D48").Value = '''7.656'
$ws.Range("E48").Value = '  +0.74%  '
$ws.Range("D49").Value = '''3.101'
$ws.Range("E49").Value = '  +3.08%  '
$ws.Range("D50").Value = '''1.844'
$ws.Range("E50").Value = '  -1.85%  '
$ws.Range("D51").Value = '2.021.31'
$ws.Range("E51").Value = '  -8.09%  '
